$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 533.6875
$ws.Range("I19").Value = 443.42856
$ws.Range("J19").Value = 603.8889
$ws.Range("K19").Value = 443.42856
$ws.Range("L19").Value = 603.8889
$ws.Range("M19").Value = -268.42856
$ws.Range("N19").Value = -953.8889

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15623.68
$ws.Range("I62").Value = 3774.625
$ws.Range("K62").Value = 3774.625
$ws.Range("M62").Value = -3150.625

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 15623.68
$ws.Range("I65").Value = 3774.625
$ws.Range("K65").Value = 18873.125
$ws.Range("M65").Value = -15753.125

# Sheet ALC, row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2500
$ws.Range("I94").Value = 2500
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2500
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2049

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4211.1665
$ws.Range("I116").Value = 2599.7144
$ws.Range("J116").Value = 9851.25
$ws.Range("K116").Value = 2599.7144
$ws.Range("L116").Value = 9851.25
$ws.Range("M116").Value = 842.2856000000002
$ws.Range("N116").Value = -16735.25

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4419.2583
$ws.Range("I132").Value = 2990.982
$ws.Range("J132").Value = 15641.429
$ws.Range("K132").Value = 8972.946
$ws.Range("L132").Value = 46924.287
$ws.Range("M132").Value = -6442.946
$ws.Range("N132").Value = -51984.287

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3206.7693
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -14100

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6940
$ws.Range("I32").Value = 7418.3213
$ws.Range("J32").Value = 1582.8
$ws.Range("K32").Value = 7418.3213
$ws.Range("L32").Value = 1582.8
$ws.Range("M32").Value = -7131.3213
$ws.Range("N32").Value = -2156.8

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3629.5789
$ws.Range("I61").Value = 2244.5
$ws.Range("J61").Value = 6004
$ws.Range("K61").Value = 2244.5
$ws.Range("L61").Value = 6004
$ws.Range("M61").Value = -2032.5
$ws.Range("N61").Value = -6428

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2496.4
$ws.Range("I74").Value = 3018.25
$ws.Range("J74").Value = 1900
$ws.Range("K74").Value = 3018.25
$ws.Range("L74").Value = 1900
$ws.Range("M74").Value = -2144.25
$ws.Range("N74").Value = -3648

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2496.4
$ws.Range("I77").Value = 3018.25
$ws.Range("J77").Value = 1900
$ws.Range("K77").Value = 15091.25
$ws.Range("L77").Value = 9500
$ws.Range("M77").Value = -10723.25
$ws.Range("N77").Value = -18236

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2511.7144
$ws.Range("I122").Value = 1637.3
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 4911.9
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -2461.9
$ws.Range("N122").Value = -64900

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4427.8887
$ws.Range("I132").Value = 1910.4286
$ws.Range("J132").Value = 6630.6665
$ws.Range("K132").Value = 5731.2858
$ws.Range("L132").Value = 19891.9995
$ws.Range("M132").Value = -3201.2858
$ws.Range("N132").Value = -24951.9995

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3629.5789
$ws.Range("I136").Value = 2244.5
$ws.Range("J136").Value = 6004
$ws.Range("K136").Value = 6733.5
$ws.Range("L136").Value = 18012
$ws.Range("M136").Value = -4183.5
$ws.Range("N136").Value = -23112

# Sheet BSM, row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 176.47058
$ws.Range("I80").Value = 170.2
$ws.Range("J80").Value = 185.42857
$ws.Range("K80").Value = 170.2
$ws.Range("L80").Value = 185.42857
$ws.Range("M80").Value = 827.8
$ws.Range("N80").Value = -2181.42857

# Sheet BSM, row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 176.47058
$ws.Range("I83").Value = 170.2
$ws.Range("J83").Value = 185.42857
$ws.Range("K83").Value = 851
$ws.Range("L83").Value = 927.1428500000001
$ws.Range("M83").Value = 4141
$ws.Range("N83").Value = -10911.14285

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1755.5
$ws.Range("I107").Value = 1755.5
$ws.Range("K107").Value = 1755.5
$ws.Range("M107").Value = 164.5

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7553.8276
$ws.Range("I134").Value = 2852.9443
$ws.Range("J134").Value = 15246.182
$ws.Range("K134").Value = 8558.832900000001
$ws.Range("L134").Value = 45738.546
$ws.Range("M134").Value = -6023.832900000001
$ws.Range("N134").Value = -50808.546

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2889.513
$ws.Range("I31").Value = 2122.5
$ws.Range("J31").Value = 3696.8948
$ws.Range("K31").Value = 2122.5
$ws.Range("L31").Value = 3696.8948
$ws.Range("M31").Value = -1827.5
$ws.Range("N31").Value = -4286.8948

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2889.513
$ws.Range("I34").Value = 2122.5
$ws.Range("J34").Value = 3696.8948
$ws.Range("K34").Value = 2122.5
$ws.Range("L34").Value = 3696.8948
$ws.Range("M34").Value = -1920.5
$ws.Range("N34").Value = -4100.8948

# Sheet CRP, row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 35500
$ws.Range("J116").Value = 35500
$ws.Range("L116").Value = 35500
$ws.Range("N116").Value = -44678

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5097.5557
$ws.Range("I132").Value = 4293.3335
$ws.Range("J132").Value = 5499.6665
$ws.Range("K132").Value = 12880.0005
$ws.Range("L132").Value = 16498.9995
$ws.Range("M132").Value = -10350.0005
$ws.Range("N132").Value = -21558.9995

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2265.5278
$ws.Range("I134").Value = 911.4091
$ws.Range("J134").Value = 4393.4287
$ws.Range("K134").Value = 2734.2273
$ws.Range("L134").Value = 13180.2861
$ws.Range("M134").Value = -199.2273
$ws.Range("N134").Value = -18250.2861

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1130.1333
$ws.Range("I5").Value = 279.27274
$ws.Range("K5").Value = 837.81822
$ws.Range("M5").Value = -725.81822

# Sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1063
$ws.Range("I86").Value = 493
$ws.Range("J86").Value = 1633
$ws.Range("K86").Value = 1479
$ws.Range("M86").Value = -293
$ws.Range("N86").Value = -7271

# Sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1063
$ws.Range("I89").Value = 493
$ws.Range("J89").Value = 1633
$ws.Range("K89").Value = 4437
$ws.Range("M89").Value = 1491
$ws.Range("N89").Value = -26553

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3392.3132
$ws.Range("I122").Value = 471.4
$ws.Range("K122").Value = 4242.599999999999
$ws.Range("M122").Value = -1792.599999999999

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1130.1333
$ws.Range("I135").Value = 279.27274
$ws.Range("K135").Value = 2513.45466
$ws.Range("M135").Value = 21.54534000000012

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2799.9524
$ws.Range("I126").Value = 3028.2307
$ws.Range("J126").Value = 2429
$ws.Range("K126").Value = 9084.6921
$ws.Range("L126").Value = 7287
$ws.Range("M126").Value = -6614.6921
$ws.Range("N126").Value = -12227

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2788.818
$ws.Range("I132").Value = 2684.1667
$ws.Range("J132").Value = 2914.4
$ws.Range("K132").Value = 8052.500100000001
$ws.Range("L132").Value = 8743.200000000001
$ws.Range("M132").Value = -5522.500100000001
$ws.Range("N132").Value = -13803.2

# Sheet LTW, row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 27415
$ws.Range("J94").Value = 27415
$ws.Range("L94").Value = 27415
$ws.Range("N94").Value = -28767

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 34490816
$ws.Range("I132").Value = 58834890
$ws.Range("J132").Value = 3387.3333
$ws.Range("K132").Value = 176504670
$ws.Range("L132").Value = 10161.9999
$ws.Range("M132").Value = -176502140
$ws.Range("N132").Value = -15221.9999

# Sheet WVR, row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6505.2
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6505.2
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6505.2
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -7487.2

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 284
$ws.Range("I107").Value = 269.82352
$ws.Range("J107").Value = 364.33334
$ws.Range("K107").Value = 809.47056
$ws.Range("L107").Value = 1093.00002
$ws.Range("M107").Value = 1110.52944
$ws.Range("N107").Value = -4933.000019999999

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3727.5715
$ws.Range("I132").Value = 2798.6667
$ws.Range("K132").Value = 8396.000100000001
$ws.Range("M132").Value = -5866.000100000001
